$wb = $excel.ActiveWorkbook

# --- Update status text "Ready for handoff" -> "In Translation" ---
# This shared string is referenced from the Overview sheet (columns for
# each locale) as well as from each per-locale detail sheet's Status column.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"

# --- Narrow the per-locale status columns ---
# Overview sheet: columns E (zh-cn) and F (de-de)
$wsOverview.Range("E1").EntireColumn.ColumnWidth = 12.5
$wsOverview.Range("F1").EntireColumn.ColumnWidth = 12.5

# zh-cn / de-de detail sheets: column C (Status)
$wsZhCn.Range("C1").EntireColumn.ColumnWidth = 12.5
$wsDeDe.Range("C1").EntireColumn.ColumnWidth = 12.5
